# Applies the "Updated cryptos list" data refresh: new Price (D) and
# Volume(1h) (E) values for rows 2-51 of the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.698.13"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "'2.094.70"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'343.14"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.5164"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "'0.09237"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "'1.162"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").Value = "'24.85"
$ws.Range("D13").Value = "'2.099.04"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "'8.284"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "'6.736"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "'99.30"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "'0.00001149"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'20.73"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "'0.06651"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").Value = "'29.733.68"
$ws.Range("E23").Value = "  -3.05%  "
$ws.Range("E24").Value = "  -3.52%  "
$ws.Range("D25").Value = "'2.321"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").Value = "'2.344.16"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'21.93"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("D28").Value = "'2.515"
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("D29").Value = "'161.24"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").Value = "'1.133"
$ws.Range("E31").Value = "  -7.62%  "
$ws.Range("D32").Value = "'0.1050"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "'1.650"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").Value = "'6.151"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").Value = "'3.938"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").Value = "'6.232"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "'10.20"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").Value = "'0.02573"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "'0.06698"
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("D40").Value = "'12.45"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").Value = "'0.6868"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").Value = "'1.321"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").Value = "'0.2221"
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("D44").Value = "'0.6689"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").Value = "'14.28"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "'0.00000000357"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").Value = "'3.616"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").Value = "'81.81"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("E51").Value = "  -1.98%  "
